$d = $word.ActiveDocument

$pairs = @(
    @("41×47=", "45×76="),
    @("51×95=", "28×74="),
    @("78×19=", "52×81="),
    @("62×69=", "58×84="),
    @("18×12=", "86×76="),
    @("30×91=", "59×33="),
    @("39×71=", "61×94="),
    @("73×99=", "71×87="),
    @("48×74=", "90×95="),
    @("25×44=", "34×53="),
    @("14×94=", "63×79="),
    @("71×82=", "25×62="),
    @("63×64=", "14×59="),
    @("69×57=", "21×69="),
    @("19×94=", "51×29="),
    @("90×13=", "69×38="),
    @("79×33=", "58×91="),
    @("71×41=", "18×17="),
    @("15×81=", "53×60="),
    @("42×98=", "50×29="),
    @("85×17=", "72×82="),
    @("44×67=", "34×74="),
    @("89×66=", "90×53="),
    @("51×68=", "66×38="),
    @("64×72=", "40×61=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
